$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug Log")
$ws.Activate()

# Copy the formatting of the row above (row 27) into the new row (28)
# so that styles/borders/number-formats match the rest of the table.
$ws.Range("A27:H27").Copy()
$ws.Range("A28:H28").PasteSpecial(-4122) # xlPasteFormats

# Fill in the new bug-log entry (bug #26)
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = 3
$ws.Cells.Item(28, 3).Value = "Boostrap "
$ws.Cells.Item(28, 4).Value = "Does not show success message or how many lines have been processed"
$ws.Cells.Item(28, 5).Value = "Resolved"
$ws.Cells.Item(28, 6).Value = "14/11/2019"
$ws.Cells.Item(28, 7).Value = "14/11/2019"
$ws.Cells.Item(28, 8).Value = "Matthew & DaEun"

# Row 28 now holds a single-line entry, so it is shorter than row 27
$ws.Rows.Item(28).RowHeight = 15.75

# Update the active selection to the newly edited cell
$ws.Range("D28").Select()
